$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (borders, fill, alignment/wrap) from row 5 down to the new row 6
$ws.Range("B5:F5").Copy()
$ws.Range("B6:F6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(5).RowHeight

# Fill in the new "BusquedaFallida" scenario row (order matters for shared-string layout)
$ws.Range("B6").Value = "BusquedaFallida"
$ws.Range("C5").Value = "En el portal de myshopify en el home utilizo el buscador y me aparece un resultado relacionado"
$ws.Range("C4").Value = "En el portal de myShopify me intento registrar sin embargo el capcha no me deja debido a que es un robot"
$ws.Range("C6").Value = "En el portal de myshopify utlizo el buscador con algo que no se encuentre en la tienda, debe aparecer el mensaje indicando que no hubieron resultados."
$ws.Range("F6").Value = "Me aparece el mensaje 'No results found for'"
$ws.Range("E6").Value = "1. ingresar en el buscador  una palabra de algo que no se vende en la tienda. 2. verificar que se encuentre el mensaje que le avisa al usuario que no se encontraro resultados."
$ws.Range("D6").Value = "Me encuentro en el home "

# Move the active selection like the author left it
$ws.Range("C10").Select()
